# Add season-record columns (Wins, Losses, Ties) to the roster worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy the formatting of the existing header cell (AC1)
# onto the three new header cells, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows (2 through 59) get the team's season record values.
for ($row = 2; $row -le 59; $row++) {
    $ws.Cells.Item($row, 30).Value = 71
    $ws.Cells.Item($row, 31).Value = 91
    $ws.Cells.Item($row, 32).Value = 0
}
